$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain plain text so numeric-looking
# strings (e.g. "572.01", "0.999") are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '62.845.80'
$ws.Range("E2").Value = '  +5.29%  '
$ws.Range("D3").Value = '3.351.56'
$ws.Range("E3").Value = '  +5.03%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '572.01'
$ws.Range("E5").Value = '  +7.00%  '
$ws.Range("D6").Value = '152.52'
$ws.Range("E6").Value = '  +5.51%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.352.75'
$ws.Range("E8").Value = '  +4.79%  '
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '7.45'
$ws.Range("E10").Value = '  +1.69%  '
$ws.Range("E11").Value = '  +5.10%  '
$ws.Range("E12").Value = '  +2.94%  '
$ws.Range("D13").Value = '3.930.23'
$ws.Range("E13").Value = '  +4.88%  '
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '26.93'
$ws.Range("E15").Value = '  +4.09%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  +4.49%  '
$ws.Range("D17").Value = '62.860.49'
$ws.Range("E17").Value = '  +5.15%  '
$ws.Range("D18").Value = '3.357.15'
$ws.Range("E18").Value = '  +5.39%  '
$ws.Range("E19").Value = '  +1.94%  '
$ws.Range("D20").Value = '13.83'
$ws.Range("E20").Value = '  +5.57%  '
$ws.Range("D21").Value = '8.41'
$ws.Range("E21").Value = '  +2.42%  '
$ws.Range("D22").Value = '384.06'
$ws.Range("E22").Value = '  +4.48%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '0.535'
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("D25").Value = '70.32'
$ws.Range("D26").Value = '9.37'
$ws.Range("E26").Value = '  +6.37%  '
$ws.Range("E27").Value = '  +6.73%  '
$ws.Range("D28").Value = '0.0₃0965'
$ws.Range("E28").Value = '  +8.73%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  +6.96%  '
$ws.Range("E31").Value = '  +11.36%  '
$ws.Range("D32").Value = '5.62'
$ws.Range("E32").Value = '  +5.90%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '23.02'
$ws.Range("E33").Value = '  +3.56%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '6.36'
$ws.Range("E34").Value = '  +4.42%  '
$ws.Range("D35").Value = '6.73'
$ws.Range("E35").Value = '  +2.47%  '
$ws.Range("D36").Value = '1.48'
$ws.Range("E36").Value = '  +11.01%  '
$ws.Range("D37").Value = '158.86'
$ws.Range("E38").Value = '  +12.57%  '
$ws.Range("D39").Value = '27.06'
$ws.Range("E39").Value = '  +5.21%  '
$ws.Range("D40").Value = '0.0329'
$ws.Range("E40").Value = '  +12.14%  '
$ws.Range("D41").Value = '0.0737'
$ws.Range("E41").Value = '  +5.62%  '
$ws.Range("D42").Value = '2.780.97'
$ws.Range("E42").Value = '  -0.49%  '
$ws.Range("D43").Value = '41.05'
$ws.Range("E43").Value = '  +4.69%  '
$ws.Range("D44").Value = '4.29'
$ws.Range("E44").Value = '  +1.82%  '
$ws.Range("D45").Value = '0.744'
$ws.Range("E45").Value = '  +4.64%  '
$ws.Range("E46").Value = '  +5.84%  '
$ws.Range("D47").Value = '3.395.45'
$ws.Range("E47").Value = '  +4.96%  '
$ws.Range("D48").Value = '21.98'
$ws.Range("E48").Value = '  +7.80%  '
$ws.Range("D49").Value = '6.34'
$ws.Range("E49").Value = '  +3.64%  '
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("D51").Value = '289.22'

# Restore default (General) formatting now that the text values are stored,
# so cell styling matches the original (unformatted) cells.
$ws.Range("D2:E51").ClearFormats()
